$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.102.88"
$ws.Range("E2").Value = "  -5.53%  "
$ws.Range("D3").Value = "2.994.36"
$ws.Range("E3").Value = "  -6.02%  "
$ws.Range("D5").Value = "'573.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.39%  "
$ws.Range("D6").Value = "'125.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.71%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "2.988.93"
$ws.Range("E8").Value = "  -6.10%  "
$ws.Range("D9").Value = "'0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "'0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.75%  "
$ws.Range("D11").Value = "'5.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.95%  "
$ws.Range("E12").Value = "  -4.04%  "
$ws.Range("D13").Value = "'0.0000220"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.80%  "
$ws.Range("D14").Value = "'32.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.42%  "
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "3.487.45"
$ws.Range("E16").Value = "  -6.02%  "
$ws.Range("D17").Value = "2.985.49"
$ws.Range("E17").Value = "  -6.34%  "
$ws.Range("D18").Value = "60.069.47"
$ws.Range("E18").Value = "  -5.65%  "
$ws.Range("D19").Value = "'6.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").Value = "'427.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.94%  "
$ws.Range("D21").Value = "'13.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.52%  "
$ws.Range("E22").Value = "  -4.79%  "
$ws.Range("D23").Value = "'7.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.49%  "
$ws.Range("D24").Value = "'12.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("D25").Value = "'79.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.11%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'2.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.78%  "
$ws.Range("D29").Value = "'1.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.47%  "
$ws.Range("D30").Value = "'7.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.32%  "
$ws.Range("E31").Value = "  -11.01%  "
$ws.Range("D32").Value = "'25.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.60%  "
$ws.Range("D33").Value = "'0.0942"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.58%  "
$ws.Range("D34").Value = "'5.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.31%  "
$ws.Range("D35").Value = "'0.931"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.18%  "
$ws.Range("D36").Value = "'50.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("E37").Value = "  -16.40%  "
$ws.Range("D38").Value = "0.0₃0664"
$ws.Range("E38").Value = "  -9.97%  "
$ws.Range("D39").Value = "'8.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("D40").Value = "'0.0355"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.03%  "
$ws.Range("E41").Value = "  -5.60%  "
$ws.Range("D42").Value = "'373.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.42%  "
$ws.Range("D43").Value = "2.674.28"
$ws.Range("E43").Value = "  -4.60%  "
$ws.Range("D44").Value = "'2.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.15%  "
$ws.Range("D45").Value = "'0.998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  -7.88%  "
$ws.Range("D47").Value = "'119.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.99%  "
$ws.Range("E48").Value = "  -7.31%  "
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("D50").Value = "'23.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.73%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.132"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.90%  "
